$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "52.269.33"
$ws.Range("E2").Value = "  +5.86%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.793.53"
$ws.Range("E3").Value = "  +6.36%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "116.78"
$ws.Range("E5").Value = "  +4.79%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "341.35"
$ws.Range("E6").Value = "  +4.92%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.553"
$ws.Range("E7").Value = "  +5.47%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("E9").Value = "  +6.14%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.03"
$ws.Range("E10").Value = "  +6.73%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0867"
$ws.Range("E11").Value = "  +7.08%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "20.13"
$ws.Range("E12").Value = "  +0.53%  "
$ws.Range("E13").Value = "  +2.73%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.64"
$ws.Range("E14").Value = "  +1.16%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.235.87"
$ws.Range("E15").Value = "  +6.19%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.809.43"
$ws.Range("E16").Value = "  +6.63%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.887"
$ws.Range("E17").Value = "  +4.21%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "52.124.15"
$ws.Range("E18").Value = "  +5.49%  "
$ws.Range("E19").Value = "  +11.03%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.34"
$ws.Range("E20").Value = "  +0.74%  "
$ws.Range("E21").Value = "  +4.56%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0₃0986"
$ws.Range("E22").Value = "  +4.27%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "278.68"
$ws.Range("E23").Value = "  +4.03%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "70.37"
$ws.Range("E24").Value = "  +2.10%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.80"
$ws.Range("E25").Value = "  +9.34%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.88"
$ws.Range("E26").Value = "  +3.49%  "
$ws.Range("E27").Value = "  -0.03%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.25"
$ws.Range("E28").Value = "  +0.97%  "
$ws.Range("E29").Value = "  +1.21%  "
$ws.Range("E30").Value = "  +3.47%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "34.87"
$ws.Range("E31").Value = "  +1.22%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "50.38"
$ws.Range("E32").Value = "  +1.75%  "
$ws.Range("E33").Value = "  +5.40%  "
$ws.Range("E34").Value = "  +2.89%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.12"
$ws.Range("E35").Value = "  +5.10%  "
$ws.Range("E36").Value = "  -0.21%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "19.00"
$ws.Range("E37").Value = "  +0.19%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.00"
$ws.Range("E38").Value = "  +1.58%  "
$ws.Range("E39").Value = "  +5.16%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0378"
$ws.Range("E40").Value = "  +13.36%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.74"
$ws.Range("E41").Value = "  +28.72%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "23.48"
$ws.Range("E42").Value = "  +4.03%  "
$ws.Range("E43").Value = "  +4.23%  "
$ws.Range("E44").Value = "  +3.17%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "124.88"
$ws.Range("E45").Value = "  -3.12%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.101.09"
$ws.Range("E46").Value = "  +2.10%  "
$ws.Range("E47").Value = "  +2.58%  "
$ws.Range("E48").Value = "  +3.53%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.58"
$ws.Range("E49").Value = "  +7.47%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.906"
$ws.Range("E50").Value = "  +21.85%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "9.00"
